$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet grows from 52 to 53 data rows (A1:E52 -> A1:E53). Every existing
# row also shifts its values (an extra historical data point was inserted at
# the top of the forecast-vector series and the y_0_forecast/y_1_forecast
# columns were recomputed), and a brand-new row 53 is appended at the bottom.
#
# Column A carries a custom date style (bold/centered/bordered, "YYYY-MM-DD
# HH:MM:SS" number format) that every existing A2:A52 cell already has. Grab
# that exact style (via copy / paste-formats-only) for the new A53 cell before
# writing any values, so the new row reuses the existing style index instead
# of Excel minting a near-duplicate one.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# date_of_forecast (A), y_0 (B), y_0_forecast (C), y_1 (D), y_1_forecast (E)
$data = @(
  @(39400, 2007, 3.145939949069287, 2008, 2.455972543253826),
  @(39583, 2008, 2.740959689118805, 2009, 3.206168778303486),
  @(39765, 2008, 1.769627576887389, 2009, 1.087227286828241),
  @(39948, 2009, -5.478010998490157, 2010, -2.013762956649334),
  @(40130, 2009, -4.774178217057779, 2010, -0.8523446516643385),
  @(40310, 2010, -0.1091898317121864, 2011, -1.305195642355683),
  @(40494, 2010, 1.97975191822708, 2011, 1.39052144387346),
  @(40676, 2011, 3.371423250978856, 2012, 0.8060632160631576),
  @(40862, 2011, 3.452886745653183, 2012, 2.798447799311043),
  @(41044, 2012, 1.627570629117536, 2013, 2.766358213445708),
  @(41228, 2012, 1.239479831392853, 2013, 2.047133666472267),
  @(41409, 2013, -0.03183655677961861, 2014, 1.102200073559878),
  @(41592, 2013, 0.2379616621361214, 2014, 1.062273708599726),
  @(41774, 2014, 1.812248956008733, 2015, 1.209672013646323),
  @(41957, 2014, 1.51977456621637, 2015, 0.6176326357196116),
  @(42137, 2015, 1.290465392296114, 2016, 0.9879295308886871),
  @(42321, 2015, 1.470039379455756, 2016, 1.577608035818301),
  @(42503, 2016, 1.57569012346459, 2017, 1.643656926428561),
  @(42689, 2016, 1.638797242243251, 2017, 1.369334405341616),
  @(42867, 2017, 1.73823635068906, 2018, 1.765380623247137),
  @(43053, 2017, 2.161565493242668, 2018, 2.486299099038347),
  @(43145, 2018, 2.50728418643813, 2019, 2.11231490846715),
  @(43235, 2018, 2.337818484846466, 2019, 2.076648015684435),
  @(43326, 2018, 2.354760705778203, 2019, 2.107524645430914),
  @(43418, 2018, 2.214251681313772, 2019, 1.362030665126834),
  @(43510, 2019, 0.8037559998091082, 2020, 1.649904670037805),
  @(43600, 2019, 0.8311911554373719, 2020, 1.758956425699298),
  @(43691, 2019, 0.7024402883234027, 2020, 1.302300993836147),
  @(43783, 2019, 0.6066442151010376, 2020, 0.7878236429522678),
  @(43875, 2020, 0.5367417164559685, 2021, 0.9721240557711397),
  @(43966, 2020, -1.538034740964334, 2021, -0.7351085756681308),
  @(44068, 2020, -4.65090747647452, 2021, -1.895157449498863),
  @(44159, 2020, -4.207901339433196, 2021, -0.4850133725290084),
  @(44251, 2021, 0.003696830084365388, 2022, -1.242205446257827),
  @(44341, 2021, 0.303920243687994, 2022, -0.6038293380915438),
  @(44432, 2021, 0.7583924418458787, 2022, 1.108416787477773),
  @(44525, 2021, 1.099928004397532, 2022, 1.102608990163567),
  @(44617, 2022, 1.880148611648913, 2023, 0.01066574587431646),
  @(44706, 2022, 1.91914784107321, 2023, 0.2384815980940092),
  @(44798, 2022, 2.236860175919531, 2023, 0.8831516962375607),
  @(44890, 2022, 2.310042359896225, 2023, 2.086256540666986),
  @(44981, 2023, 0.2425620590337463, 2024, 1.106935253696562),
  @(45071, 2023, -0.1211988132392205, 2024, 0.7707080878861294),
  @(45163, 2023, -0.09588622947416248, 2024, 0.8742015250004842),
  @(45254, 2023, 0.0464415346324687, 2024, 0.3491198177708599),
  @(45345, 2024, -0.297474409307219, 2025, 0.07171493608653101),
  @(45436, 2024, -0.04760886976447054, 2025, 0.5568966348730831),
  @(45534, 2024, -0.2385784141923808, 2025, -0.06923086958923186),
  @(45618, 2024, -0.3101476031197148, 2025, 0.2126457877301924),
  @(45713, 2025, -0.07765238411295838, 2026, -0.17084471036517),
  @(45800, 2025, 0.1245593350339691, 2026, 0.1722027100061974),
  @(45891, 2025, -0.0960403240804597, 2026, -0.1059622177528863)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r++
}

Write-Host "Updated rows 2-53 with corrected forecast data"
